$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 951.8182
$ws.Range("J17").Value = 951.8182
$ws.Range("L17").Value = 2855.4546
$ws.Range("N17").Value = -3191.4546
$ws.Range("H33").Value = 250.21428
$ws.Range("I33").Value = 238.6923
$ws.Range("K33").Value = 238.6923
$ws.Range("M33").Value = -9.692299999999989
$ws.Range("H51").Value = 6121.75
$ws.Range("I51").Value = 4495.3335
$ws.Range("K51").Value = 4495.3335
$ws.Range("M51").Value = -4011.3335
$ws.Range("H58").Value = 6337.643
$ws.Range("J58").Value = 9533.666999999999
$ws.Range("L58").Value = 28601.001
$ws.Range("N58").Value = -28901.001
$ws.Range("H100").Value = 7544
$ws.Range("I100").Value = 1815.9
$ws.Range("J100").Value = 19000.2
$ws.Range("K100").Value = 1815.9
$ws.Range("L100").Value = 19000.2
$ws.Range("M100").Value = -1274.9
$ws.Range("N100").Value = -20082.2
$ws.Range("H132").Value = 3374.8572
$ws.Range("I132").Value = 2024.8
$ws.Range("K132").Value = 6074.4
$ws.Range("M132").Value = -3544.4
$ws.Range("H135").Value = 2296.1428
$ws.Range("I135").Value = 954.05884
$ws.Range("K135").Value = 8586.529560000001
$ws.Range("M135").Value = -6051.529560000001
$ws.Range("H138").Value = 5367.0586
$ws.Range("J138").Value = 6385.1714
$ws.Range("L138").Value = 19155.5142
$ws.Range("N138").Value = -29435.5142
$ws.Range("H141").Value = 6267.5
$ws.Range("I141").Value = 5995
$ws.Range("J141").Value = 6358.3335
$ws.Range("K141").Value = 17985
$ws.Range("L141").Value = 19075.0005
$ws.Range("M141").Value = -12805
$ws.Range("N141").Value = -29435.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3681
$ws.Range("H61").Value = 2416.3845
$ws.Range("I61").Value = 1786.5714
$ws.Range("K61").Value = 1786.5714
$ws.Range("M61").Value = -1574.5714
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H136").Value = 2416.3845
$ws.Range("I136").Value = 1786.5714
$ws.Range("K136").Value = 5359.7142
$ws.Range("M136").Value = -2809.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1200.5
$ws.Range("I86").Value = 699.5
$ws.Range("J86").Value = 2703.5
$ws.Range("K86").Value = 699.5
$ws.Range("L86").Value = 2703.5
$ws.Range("M86").Value = 423.5
$ws.Range("N86").Value = -4949.5
$ws.Range("H89").Value = 1200.5
$ws.Range("I89").Value = 699.5
$ws.Range("J89").Value = 2703.5
$ws.Range("K89").Value = 3497.5
$ws.Range("L89").Value = 13517.5
$ws.Range("M89").Value = 2118.5
$ws.Range("N89").Value = -24749.5
$ws.Range("H99").Value = 4279.2
$ws.Range("I99").Value = 3799.3333
$ws.Range("K99").Value = 3799.3333
$ws.Range("M99").Value = -2301.3333
$ws.Range("H134").Value = 4199.5
$ws.Range("J134").Value = 7639
$ws.Range("L134").Value = 22917
$ws.Range("N134").Value = -27987

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3508.4285
$ws.Range("I31").Value = 1352.875
$ws.Range("K31").Value = 1352.875
$ws.Range("M31").Value = -1057.875
$ws.Range("H34").Value = 3508.4285
$ws.Range("I34").Value = 1352.875
$ws.Range("K34").Value = 1352.875
$ws.Range("M34").Value = -1150.875
$ws.Range("H58").Value = 336360.44
$ws.Range("I58").Value = 527457.6
$ws.Range("K58").Value = 527457.6
$ws.Range("M58").Value = -527254.6
$ws.Range("H103").Value = 7524.8
$ws.Range("I103").Value = 7524.8
$ws.Range("K103").Value = 7524.8
$ws.Range("M103").Value = -6352.8
$ws.Range("H132").Value = 4152.381
$ws.Range("J132").Value = 5321.1665
$ws.Range("L132").Value = 15963.4995
$ws.Range("N132").Value = -21023.4995
$ws.Range("H136").Value = 336360.44
$ws.Range("I136").Value = 527457.6
$ws.Range("K136").Value = 1582372.8
$ws.Range("M136").Value = -1579822.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1905.6666
$ws.Range("I5").Value = 1998
$ws.Range("J5").Value = 1859.5
$ws.Range("K5").Value = 5994
$ws.Range("L5").Value = 5578.5
$ws.Range("M5").Value = -5882
$ws.Range("N5").Value = -5802.5
$ws.Range("H12").Value = 282.6875
$ws.Range("J12").Value = 312.57144
$ws.Range("L12").Value = 937.71432
$ws.Range("N12").Value = -1283.71432
$ws.Range("H113").Value = 5292170
$ws.Range("I113").Value = 12346112
$ws.Range("J113").Value = 1713.25
$ws.Range("K113").Value = 37038336
$ws.Range("L113").Value = 5139.75
$ws.Range("M113").Value = -37036166
$ws.Range("N113").Value = -9479.75
$ws.Range("H131").Value = 3575.795
$ws.Range("J131").Value = 3807.647
$ws.Range("L131").Value = 11422.941
$ws.Range("N131").Value = -21502.941
$ws.Range("H135").Value = 1905.6666
$ws.Range("I135").Value = 1998
$ws.Range("J135").Value = 1859.5
$ws.Range("K135").Value = 17982
$ws.Range("L135").Value = 16735.5
$ws.Range("M135").Value = -15447
$ws.Range("N135").Value = -21805.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3002
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 3002
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3002
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -3226
$ws.Range("H63").Value = 16701
$ws.Range("J63").Value = 16701
$ws.Range("L63").Value = 16701
$ws.Range("N63").Value = -18073
$ws.Range("H66").Value = 16701
$ws.Range("J66").Value = 16701
$ws.Range("L66").Value = 50103
$ws.Range("N66").Value = -56967
$ws.Range("H80").Value = 1431763.9
$ws.Range("I80").Value = 1431153
$ws.Range("J80").Value = 1432374.8
$ws.Range("K80").Value = 1431153
$ws.Range("L80").Value = 1432374.8
$ws.Range("M80").Value = -1430155
$ws.Range("N80").Value = -1434370.8
$ws.Range("H83").Value = 1431763.9
$ws.Range("I83").Value = 1431153
$ws.Range("J83").Value = 1432374.8
$ws.Range("K83").Value = 7155765
$ws.Range("L83").Value = 7161874
$ws.Range("M83").Value = -7150773
$ws.Range("N83").Value = -7171858
$ws.Range("H113").Value = 636656.5
$ws.Range("I113").Value = 773962.25
$ws.Range("J113").Value = 41665
$ws.Range("K113").Value = 773962.25
$ws.Range("L113").Value = 41665
$ws.Range("M113").Value = -771792.25
$ws.Range("N113").Value = -46005
$ws.Range("H126").Value = 250002320
$ws.Range("J126").Value = 6500
$ws.Range("L126").Value = 19500
$ws.Range("N126").Value = -24440
$ws.Range("H132").Value = 246563.83
$ws.Range("I132").Value = 259047.56
$ws.Range("J132").Value = 3131
$ws.Range("K132").Value = 777142.6799999999
$ws.Range("L132").Value = 9393
$ws.Range("M132").Value = -774612.6799999999
$ws.Range("N132").Value = -14453

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 737.4
$ws.Range("I22").Value = 731.2857
$ws.Range("J22").Value = 751.6667
$ws.Range("K22").Value = 731.2857
$ws.Range("L22").Value = 751.6667
$ws.Range("M22").Value = -436.2857
$ws.Range("N22").Value = -1341.6667
$ws.Range("H27").Value = 737.4
$ws.Range("I27").Value = 731.2857
$ws.Range("J27").Value = 751.6667
$ws.Range("K27").Value = 731.2857
$ws.Range("L27").Value = 751.6667
$ws.Range("M27").Value = -624.2857
$ws.Range("N27").Value = -965.6667
$ws.Range("H68").Value = 114684.555
$ws.Range("I68").Value = 4020.125
$ws.Range("K68").Value = 4020.125
$ws.Range("M68").Value = -3271.125
$ws.Range("H71").Value = 114684.555
$ws.Range("I71").Value = 4020.125
$ws.Range("K71").Value = 20100.625
$ws.Range("M71").Value = -16356.625
$ws.Range("H132").Value = 5760.615
$ws.Range("I132").Value = 4482.1665
$ws.Range("K132").Value = 13446.4995
$ws.Range("M132").Value = -10916.4995
$ws.Range("H136").Value = 3243.4055
$ws.Range("I136").Value = 3156.4783
$ws.Range("K136").Value = 9469.4349
$ws.Range("M136").Value = -6919.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 17333332
$ws.Range("J5").Value = 17333332
$ws.Range("L5").Value = 17333332
$ws.Range("N5").Value = -17333556
$ws.Range("H125").Value = 67442.5
$ws.Range("J125").Value = 67442.5
$ws.Range("L125").Value = 67442.5
$ws.Range("N125").Value = -77282.5
$ws.Range("H128").Value = 74000
$ws.Range("J128").Value = 74000
$ws.Range("L128").Value = 74000
$ws.Range("N128").Value = -83960
$ws.Range("H131").Value = 127116
$ws.Range("J131").Value = 127116
$ws.Range("L131").Value = 127116
$ws.Range("N131").Value = -137196
$ws.Range("H132").Value = 4915.7915
$ws.Range("I132").Value = 4373.8125
$ws.Range("K132").Value = 13121.4375
$ws.Range("M132").Value = -10591.4375
